# Chronologie2023.xlsx — add 5 new observation rows, then re-sort the data
# table by Date (asc) then Zone climatique (asc), matching a manual
# "add rows at bottom, then re-apply Data > Sort" edit in Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Append 5 new rows (53-57) at the bottom of the table, seeded from
#    row 52's formatting (same body style used throughout the table),
#    then overwrite the values.
# ---------------------------------------------------------------------
$ws.Range("A52:I52").Copy($ws.Range("A53:I57"))

# Row 53: RASY @ Sainte-Anne-des-Monts (Gaspésie), zone D, Cote 3
$ws.Range("A53").Value = "2023-05-08"
$ws.Range("B53").Value = "RASY"
$ws.Range("C53").Value = 132
$ws.Range("D53").Value = "Sainte-Anne-des-Monts"
$ws.Range("E53").Value = "Gaspésie"
$ws.Range("F53").Value = "D"
$ws.Range("G53").Value = "Cote 3"
$ws.Range("H53").ClearContents()
$ws.Range("I53").Value = "Jean-Philippe Baillargeon"

# Row 54: PSCR @ Sainte-Anne-des-Monts (Gaspésie), zone D, Cote 3
$ws.Range("A54").Value = "2023-05-08"
$ws.Range("B54").Value = "PSCR"
$ws.Range("C54").Value = 132
$ws.Range("D54").Value = "Sainte-Anne-des-Monts"
$ws.Range("E54").Value = "Gaspésie"
$ws.Range("F54").Value = "D"
$ws.Range("G54").Value = "Cote 3"
$ws.Range("H54").ClearContents()
$ws.Range("I54").Value = "Jean-Philippe Baillargeon"

# Row 55: BUAM @ La Tuque (Mauricie), zone C, Cote 1
$ws.Range("A55").Value = "2023-05-10"
$ws.Range("B55").Value = "BUAM"
$ws.Range("C55").Value = "N/A"
$ws.Range("D55").Value = "La Tuque"
$ws.Range("E55").Value = "Mauricie"
$ws.Range("F55").Value = "C"
$ws.Range("G55").Value = "Cote 1"
$ws.Range("H55").Value = "iNaturalist (https://www.inaturalist.org/observations/160902825)"
$ws.Range("I55").Value = "Lyse Lafrenière"

# Row 56: BUAM @ Lévis (Chaudière-Appalaches), zone B, Cote 1
$ws.Range("A56").Value = "2023-05-12"
$ws.Range("B56").Value = "BUAM"
$ws.Range("C56").Value = "N/A"
$ws.Range("D56").Value = "Lévis"
$ws.Range("E56").Value = "Chaudière-Appalaches"
$ws.Range("F56").Value = "B"
$ws.Range("G56").Value = "Cote 1"
$ws.Range("H56").Value = "Donnée soumise à l'AARQ"
$ws.Range("I56").Value = "Jean Rodrigue"

# Row 57: RASY @ Forillon (Gaspésie), zone D, Cote 3
$ws.Range("A57").Value = "2023-05-12"
$ws.Range("B57").Value = "RASY"
$ws.Range("C57").Value = 212
$ws.Range("D57").Value = "Forillon"
$ws.Range("E57").Value = "Gaspésie"
$ws.Range("F57").Value = "D"
$ws.Range("G57").Value = "Cote 3"
$ws.Range("H57").ClearContents()
$ws.Range("I57").Value = "Diane Ostiguy"

# ---------------------------------------------------------------------
# 2. Fix up the "Espèce" (B) and "Zone climatique" (F) fill colours so
#    the new rows carry the same colour-coding as existing RASY/PSCR/BUAM
#    and A/B/C/D/E/F zone cells elsewhere in the table.
# ---------------------------------------------------------------------
$ws.Range("B4").Copy()
$ws.Range("B53").PasteSpecial(-4122)
$ws.Range("B57").PasteSpecial(-4122)

$ws.Range("B6").Copy()
$ws.Range("B54").PasteSpecial(-4122)

$ws.Range("B10").Copy()
$ws.Range("B55").PasteSpecial(-4122)
$ws.Range("B56").PasteSpecial(-4122)

$ws.Range("F11").Copy()
$ws.Range("F53").PasteSpecial(-4122)
$ws.Range("F54").PasteSpecial(-4122)
$ws.Range("F57").PasteSpecial(-4122)

$ws.Range("F13").Copy()
$ws.Range("F55").PasteSpecial(-4122)

$ws.Range("F16").Copy()
$ws.Range("F56").PasteSpecial(-4122)

# Restore values clobbered by the format-only paste (PasteSpecial formats
# shouldn't touch values, but re-assert them defensively).
$ws.Range("B53").Value = "RASY"
$ws.Range("B54").Value = "PSCR"
$ws.Range("B55").Value = "BUAM"
$ws.Range("B56").Value = "BUAM"
$ws.Range("B57").Value = "RASY"
$ws.Range("F53").Value = "D"
$ws.Range("F54").Value = "D"
$ws.Range("F55").Value = "C"
$ws.Range("F56").Value = "B"
$ws.Range("F57").Value = "D"

$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# 3. Re-sort the whole data table (A4:I57) by Date ascending, then Zone
#    climatique ascending — same as the existing sortState condition,
#    just extended to cover the 5 freshly-added rows.
# ---------------------------------------------------------------------
$sortRange = $ws.Range("A4:I57")
$keyDate = $ws.Range("A4:A57")
$keyZone = $ws.Range("F4:F57")

$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add2($keyDate, $null, 1, $null, 1)
$sortObj.SortFields.Add2($keyZone, $null, 1, $null, 1)
$sortObj.SetRange($sortRange)
$sortObj.Header = 2
$sortObj.Apply()

# ---------------------------------------------------------------------
# 4. Restore the view: scroll near the new rows and select A1:B1 (the
#    merged title cell) like the saved workbook does.
# ---------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A1:B1").Select()

Write-Output "done"
